# "Basic action -> verb flow"
# - Reword the "Build system ... matching ActionEvents" task to
#   "Build basic system ... matching ActionEvents"
# - Mark that task as done (strike-through, same yellow highlight)
# - Widen column A slightly to fit the new text
# - Move the active selection down to A9 (the next empty row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 holds "Build system for matching ActionEvents to verb/primitive
# definitions" - reword it and mark the task complete (strikethrough),
# keeping its existing yellow fill.
$cell = $ws.Range("A6")
$cell.Value2 = "Build basic system for matching ActionEvents to verb/primitive definitions"
$cell.Font.Strikethrough = $true

# Widen column A so the (now longer) task text still fits.
$ws.Columns("A").ColumnWidth = 60.15

# Move the selection to the next free row.
$ws.Range("A9").Select() | Out-Null
